$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is being extended with a new "2021" column (O), mirroring the
# existing "2020" column (N) for formatting, then the real figures are
# written in on top.

# Row 2 (the thin separator/border row under the header) - blank cell,
# same formatting as its row neighbours.
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# Row 3 - year header "2021"
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 2021

# Row 4 - per-capita figure, computed the same way column N is
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Formula = "=O5/O6*1000"

# Row 5 - total waste removed (thousand tons). Picks up the plain
# ("Items"-column-style) look rather than the neighbouring N5 look.
$ws.Range("B5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 1229.5999999999999

# Row 6 - average annual population (thousand people)
$ws.Range("N6").Copy()
$ws.Range("O6").PasteSpecial(-4122)
$ws.Range("O6").Value = 6436.9

$excel.CutCopyMode = 0

# Selection moves on (as captured by the saved view state)
$ws.Range("P16").Select() | Out-Null
